$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.404.06'
$ws.Range("D3").Value = '1.638.28'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '299.40'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -1.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3782'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -0.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3514'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.53'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -3.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08060'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.207'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -3.63%  '
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.96'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -3.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.332'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -3.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.271'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -2.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001199'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("D17").Value = '1.642.85'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.93'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06937'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.702'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.28'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -2.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.28'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -4.29%  '
$ws.Range("D24").Value = '23.413.05'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.477'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -1.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.900'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -5.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.80'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.24'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.189'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -1.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.39'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -1.24%  '
$ws.Range("D31").Value = '1.820.50'
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.803'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.138'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -4.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.39'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -3.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9749'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -8.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02685'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -4.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08724'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2419'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -4.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.855'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -4.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06781'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -4.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.88'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6830'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -3.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.308'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -2.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.44'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -3.69%  '
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6315'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -3.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.239'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -4.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.901'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07679'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -3.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '126.78'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.137'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -4.55%  '
